$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: cells whose new value is a plain number-looking string must be
# forced to Text so Excel does not auto-convert them to a numeric type,
# matching the original inline-string cell type. NumberFormat is reset to
# "@" just long enough to type the value, then the cell Style is restored to
# "Normal" so the cell keeps the default (unstyled) look, same as the source.

$ws.Range("D2").Value = '40.111.41'
$ws.Range("E2").Value = '  +0.19%  '

$ws.Range("D3").Value = '2.222.22'
$ws.Range("E3").Value = '  +0.24%  '

$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '293.92'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.16%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '87.92'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.63%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.514'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.37%  '

$ws.Range("E8").Value = '  -0.03%  '

$ws.Range("E9").Value = '  -0.41%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '30.78'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.07%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '50.53'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.44%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0782'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.12%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.114'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.53%  '

$ws.Range("E14").Value = '  -0.35%  '

$ws.Range("D15").Value = '2.581.71'
$ws.Range("E15").Value = '  +0.82%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '13.85'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.24%  '

$ws.Range("D17").Value = '2.213.00'
$ws.Range("E17").Value = '  +0.14%  '

$ws.Range("E18").Value = '  +1.06%  '

$ws.Range("D19").Value = '40.053.09'
$ws.Range("E19").Value = '  +0.17%  '

$ws.Range("E20").Value = '  +0.61%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.29'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.95%  '

$ws.Range("E22").Value = '  -0.55%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '65.74'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.17%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '236.54'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.39%  '

$ws.Range("E25").Value = '  -0.10%  '

$ws.Range("E26").Value = '  +0.88%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.84'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.21%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '23.18'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.45%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.34'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.00%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.06'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -6.49%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '158.72'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.69%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '31.84'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.33%  '

$ws.Range("E33").Value = '  -0.05%  '

$ws.Range("E34").Value = '  -0.05%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.03'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +6.37%  '

$ws.Range("E36").Value = '  -0.63%  '

$ws.Range("E37").Value = '  -2.49%  '

$ws.Range("E38").Value = '  +1.36%  '

$ws.Range("E39").Value = '  +2.62%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0994'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.69%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '15.72'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.07%  '

$ws.Range("D42").Value = '2.083.70'
$ws.Range("E42").Value = '  -0.88%  '

$ws.Range("E43").Value = '  -2.51%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '19.08'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +7.23%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.09'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.26%  '

$ws.Range("E46").Value = '  +0.62%  '

$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.95'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -10.63%  '

$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.73'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.59%  '

$ws.Range("D49").Value = '2.449.89'
$ws.Range("E49").Value = '  +0.70%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.47'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.14%  '

$ws.Range("E51").Value = '  +3.86%  '
